$d = $word.ActiveDocument

$d.Content.Find.Execute("118÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "521÷3=", 2) | Out-Null
$d.Content.Find.Execute("769÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "669÷6=", 2) | Out-Null
$d.Content.Find.Execute("941÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "708÷7=", 2) | Out-Null
$d.Content.Find.Execute("999÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "967÷4=", 2) | Out-Null
$d.Content.Find.Execute("251÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "255÷8=", 2) | Out-Null
$d.Content.Find.Execute("393÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "788÷9=", 2) | Out-Null
$d.Content.Find.Execute("977÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "853÷3=", 2) | Out-Null
$d.Content.Find.Execute("962÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "110÷2=", 2) | Out-Null
$d.Content.Find.Execute("203÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "465÷9=", 2) | Out-Null
$d.Content.Find.Execute("994÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "573÷7=", 2) | Out-Null
$d.Content.Find.Execute("218÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "755÷8=", 2) | Out-Null
$d.Content.Find.Execute("890÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "117÷3=", 2) | Out-Null
$d.Content.Find.Execute("584÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "581÷2=", 2) | Out-Null
$d.Content.Find.Execute("275÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "198÷5=", 2) | Out-Null
$d.Content.Find.Execute("728÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "101÷6=", 2) | Out-Null
$d.Content.Find.Execute("280÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "497÷6=", 2) | Out-Null
$d.Content.Find.Execute("940÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "494÷8=", 2) | Out-Null
$d.Content.Find.Execute("722÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "167÷4=", 2) | Out-Null
$d.Content.Find.Execute("820÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "305÷7=", 2) | Out-Null
$d.Content.Find.Execute("262÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "997÷6=", 2) | Out-Null
$d.Content.Find.Execute("626÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "645÷8=", 2) | Out-Null
$d.Content.Find.Execute("741÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "214÷5=", 2) | Out-Null
$d.Content.Find.Execute("681÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "249÷8=", 2) | Out-Null
$d.Content.Find.Execute("754÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "194÷2=", 2) | Out-Null
$d.Content.Find.Execute("439÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "190÷6=", 2) | Out-Null
